# Weekly update: insert a new price record as the first row of the
# "Feria Lagunitas de Puerto Montt - Cilantro" weekly log (row 107),
# pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 107, shifting rows 107:163 down to 108:164.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the latest week's record.
$ws.Range("A107").Value = 4
$ws.Range("B107").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value = "Los Lagos"
$ws.Range("D107").Value = 44460
$ws.Range("E107").Value = 10
$ws.Range("F107").Value = 100112040
$ws.Range("G107").Value = "Cilantro"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 12000
$ws.Range("L107").Value = 12000
$ws.Range("M107").Value = 12000
$ws.Range("N107").Value = "$/caja 36 atados"
$ws.Range("O107").Value = "Región Metropolitana"
$ws.Range("P107").Value = 333
$ws.Range("Q107").Value = 36
$ws.Range("R107").Value = "Hortaliza"
